$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "title"
$ws.Range("B1").Value = "description"
$ws.Range("C1").Value = "image"
